$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D and E columns for data rows (2-51) stay as Text so numeric-looking
# strings (e.g. '22.80', '0.631') are not auto-converted to numbers and lose
# formatting such as trailing zeros.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '37.315.39'
$ws.Range('E2').Value = '  +2.32%  '
$ws.Range('D3').Value = '2.005.67'
$ws.Range('E3').Value = '  +2.90%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '247.67'
$ws.Range('E5').Value = '  +1.83%  '
$ws.Range('D6').Value = '0.631'
$ws.Range('E6').Value = '  +3.02%  '
$ws.Range('D7').Value = '60.56'
$ws.Range('E7').Value = '  +4.60%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  +1.98%  '
$ws.Range('D10').Value = '0.0804'
$ws.Range('E10').Value = '  +2.13%  '
$ws.Range('E11').Value = '  +1.33%  '
$ws.Range('D12').Value = '15.07'
$ws.Range('E12').Value = '  +9.90%  '
$ws.Range('D13').Value = '22.80'
$ws.Range('E13').Value = '  +7.08%  '
$ws.Range('D14').Value = '2.299.07'
$ws.Range('E14').Value = '  +2.88%  '
$ws.Range('D15').Value = '0.848'
$ws.Range('E15').Value = '  +2.88%  '
$ws.Range('D16').Value = '5.45'
$ws.Range('E16').Value = '  +3.44%  '
$ws.Range('D17').Value = '2.006.99'
$ws.Range('E17').Value = '  +3.13%  '
$ws.Range('D18').Value = '37.267.71'
$ws.Range('E18').Value = '  +2.47%  '
$ws.Range('D19').Value = '70.54'
$ws.Range('E19').Value = '  +2.01%  '
$ws.Range('D20').Value = '0.0₃0867'
$ws.Range('E20').Value = '  +2.54%  '
$ws.Range('D21').Value = '5.20'
$ws.Range('E21').Value = '  +3.94%  '
$ws.Range('D22').Value = '231.05'
$ws.Range('E22').Value = '  +1.38%  '
$ws.Range('E23').Value = '  +0.14%  '
$ws.Range('E24').Value = '  +1.39%  '
$ws.Range('D25').Value = '2.36'
$ws.Range('E25').Value = '  +0.35%  '
$ws.Range('D26').Value = '0.145'
$ws.Range('E26').Value = '  +6.97%  '
$ws.Range('D27').Value = '9.43'
$ws.Range('E27').Value = '  +4.03%  '
$ws.Range('D28').Value = '163.98'
$ws.Range('E28').Value = '  +2.22%  '
$ws.Range('D29').Value = '19.72'
$ws.Range('E29').Value = '  +2.31%  '
$ws.Range('D30').Value = '1.33'
$ws.Range('E30').Value = '  +13.04%  '
$ws.Range('E31').Value = '  +1.54%  '
$ws.Range('D32').Value = '4.85'
$ws.Range('E32').Value = '  +4.00%  '
$ws.Range('D33').Value = '0.0653'
$ws.Range('E33').Value = '  +7.44%  '
$ws.Range('D34').Value = '4.55'
$ws.Range('E34').Value = '  +5.65%  '
$ws.Range('E35').Value = '  +6.31%  '
$ws.Range('E36').Value = '  +0.03%  '
$ws.Range('E37').Value = '  +2.50%  '
$ws.Range('D38').Value = '3.29'
$ws.Range('E38').Value = '  -4.14%  '
$ws.Range('D39').Value = '5.52'
$ws.Range('E39').Value = '  +5.19%  '
$ws.Range('E40').Value = '  +0.83%  '
$ws.Range('E41').Value = '  +1.04%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').Value = '0.0215'
$ws.Range('E42').Value = '  +2.90%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value = '1.19'
$ws.Range('E43').Value = '  +2.77%  '
$ws.Range('D44').Value = '16.70'
$ws.Range('E44').Value = '  +6.33%  '
$ws.Range('D45').Value = '90.96'
$ws.Range('E45').Value = '  +4.34%  '
$ws.Range('D46').Value = '1.372.70'
$ws.Range('E46').Value = '  +1.07%  '
$ws.Range('E47').Value = '  +3.04%  '
$ws.Range('E48').Value = '  +2.52%  '
$ws.Range('E49').Value = '  +1.05%  '
$ws.Range('E50').Value = '  +15.74%  '
$ws.Range('D51').Value = '46.14'
$ws.Range('E51').Value = '  +5.99%  '
